# Rename the existing sheet to "ValidLogin" and populate it with the
# valid-login test data (header row, a couple of valid admin/manager
# credentials, and a duplicated "trainee" row).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ValidLogin"

# Data is entered in the same order the original author appears to have
# used (credentials first, header afterwards, then the trainee row) so
# that the generated shared-string table lines up with the source file.
$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "manager"
$ws1.Range("A1").Value = "Username"
$ws1.Range("B1").Value = "Password"
$ws1.Range("A3").Value = "trainee"
$ws1.Range("B3").Value = "trainee"

$ws1.PageSetup.Orientation = 1

# Add a second sheet right after "ValidLogin" for the invalid-login
# credentials test data.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "InvalidLogin"
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "password"
$ws2.Range("A2").Value = "abcd"
$ws2.Range("B2").Value = "xyz"

# Leave "InvalidLogin" as the active sheet/selected cell, matching the
# workbook's saved view state.
$ws2.Range("B3").Select()
